$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '27.500.27'
$ws.Range("E2").Value = '  +3.72%  '

# Row 3
$ws.Range("D3").Value = '1.753.61'
$ws.Range("E3").Value = '  +1.91%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9965'
$ws.Range("E4").Value = '  +0.13%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '242.69'
$ws.Range("E5").Value = '  +0.99%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9968'
$ws.Range("E6").Value = '  +0.10%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4822'
$ws.Range("E7").Value = '  -1.64%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2625'
$ws.Range("E8").Value = '  +1.40%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06182'
$ws.Range("E9").Value = '  -0.10%  '

# Row 10
$ws.Range("D10").Value = '1.742.09'
$ws.Range("E10").Value = '  +1.20%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '16.15'
$ws.Range("E11").Value = '  +3.59%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.06939'
$ws.Range("E12").Value = '  -0.13%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6062'
$ws.Range("E13").Value = '  +0.31%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.488'
$ws.Range("E14").Value = '  +0.63%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '77.52'
$ws.Range("E15").Value = '  +1.20%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9990'
$ws.Range("E16").Value = '  +0.29%  '

# Row 17
$ws.Range("D17").Value = '27.426.88'
$ws.Range("E17").Value = '  +4.08%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9972'
$ws.Range("E18").Value = '  +0.23%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007085'
$ws.Range("E19").Value = '  -0.50%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.52'
$ws.Range("E20").Value = '  +1.86%  '

# Row 21
$ws.Range("D21").Value = '1.973.28'
$ws.Range("E21").Value = '  +1.61%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.459'
$ws.Range("E22").Value = '  +1.20%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.452'
$ws.Range("E23").Value = '  +0.33%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.110'
$ws.Range("E24").Value = '  +0.66%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '141.80'
$ws.Range("E25").Value = '  +2.75%  '

# Row 26
$ws.Range("B26").Value = 'LidoDAOToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.852'
$ws.Range("E26").Value = '  +6.35%  '

# Row 27
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.29'
$ws.Range("E27").Value = '  +0.38%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '108.83'
$ws.Range("E28").Value = '  +2.93%  '

# Row 29
$ws.Range("E29").Value = '  -0.45%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.967'
$ws.Range("E30").Value = '  +1.54%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.07997'
$ws.Range("E31").Value = '  +1.00%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.689'
$ws.Range("E32").Value = '  +1.89%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04691'
$ws.Range("E33").Value = '  +4.71%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.598'
$ws.Range("E34").Value = '  -0.21%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.016'
$ws.Range("E35").Value = '  +1.92%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6214'
$ws.Range("E36").Value = '  +0.58%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9257'
$ws.Range("E37").Value = '  -2.43%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.549'
$ws.Range("E38").Value = '  +6.62%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.008'
$ws.Range("E39").Value = '  +0.12%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.0000'
$ws.Range("E40").Value = '  +0.50%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.748'
$ws.Range("E41").Value = '  +5.16%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.01499'
$ws.Range("E42").Value = '  +0.89%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '99.94'

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3856'
$ws.Range("E44").Value = '  +1.01%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '6.913'
$ws.Range("E45").Value = '  +0.32%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1159'
$ws.Range("E46").Value = '  +0.17%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05366'

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.828'
$ws.Range("E48").Value = '  +0.94%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '29.94'
$ws.Range("E49").Value = '  -1.65%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.251'
$ws.Range("E50").Value = '  +3.30%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '51.18'
$ws.Range("E51").Value = '  -0.19%  '
